# DOMA-8525: allow to restrict send meter readings
# - add a new "Автоматический" (isAutomatic) column after "Место установки счетчика"
# - store the personal-account / meter-number / tariff-count example columns as text
#   (so values like "111" round-trip as strings, not numbers)
# - store the example transmission/verification dates as plain text values instead of
#   real date-formatted serials
# - bump one sample meter reading (100 -> 100.5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column T: "Автоматический" -----------------------------------
# Clone the formatting of column S (last existing column) onto the new column T
# so the new header/cells pick up the same fill/border style, then set the header text.
$ws.Range("S1:S11").Copy()
$ws.Range("T1:T11").PasteSpecial(-4122)
$ws.Cells.Item(1, 20).Value = "Автоматический"
$ws.Columns.Item(20).ColumnWidth = $ws.Columns.Item(19).ColumnWidth

# --- Columns B (Лицевой счет), D (Номер счетчика), F (Количество тарифов): --
# convert the numeric example values to text so they are stored as strings.
for ($r = 2; $r -le 11; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = [string]$bVal

    $dVal = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = [string]$dVal

    $fVal = $ws.Cells.Item($r, 6).Value()
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = [string]$fVal
}

# --- Column M (Дата поверки): store as plain text "2021-12-20" ------------
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 13).NumberFormat = "@"
    $ws.Cells.Item($r, 13).Value = "2021-12-20"
}

# --- L7 (Дата передачи показаний) was a real date value; make it text too -
$ws.Cells.Item(7, 12).NumberFormat = "@"
$ws.Cells.Item(7, 12).Value = "2021-12-20"

# --- H2 (Показание 1): sample reading tweaked from 100 to 100.5 -----------
$ws.Cells.Item(2, 8).Value = 100.5

# --- Columns I, J, K: normalize blank-cell number format ------------------
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 9).NumberFormat = "General"
    $ws.Cells.Item($r, 10).NumberFormat = "General"
    $ws.Cells.Item($r, 11).NumberFormat = "General"
}
